$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching style of existing header cells (B1:E1) -
# bold font, thin box border, centered/top aligned - same as the rest of row 1
$ws.Range("F1").Value = "time_taken"
$headerCell = $ws.Range("F1")
$headerCell.Font.Bold = $true
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = -4108
$headerCell.VerticalAlignment = -4160

# Timestamps recorded for each data row (rows 2-23), taken from the source diff
$timestamps = @(
    "2021-10-05 10:52:06.350550",
    "2021-10-05 10:52:06.350563",
    "2021-10-05 10:52:06.350567",
    "2021-10-05 10:52:06.350570",
    "2021-10-05 10:52:06.350574",
    "2021-10-05 10:52:06.350577",
    "2021-10-05 10:52:06.350580",
    "2021-10-05 10:52:06.350583",
    "2021-10-05 10:52:06.350587",
    "2021-10-05 10:52:06.350590",
    "2021-10-05 10:52:06.350593",
    "2021-10-05 10:52:06.350596",
    "2021-10-05 10:52:06.350599",
    "2021-10-05 10:52:06.350602",
    "2021-10-05 10:52:06.350605",
    "2021-10-05 10:52:06.350608",
    "2021-10-05 10:52:06.350611",
    "2021-10-05 10:52:06.350614",
    "2021-10-05 10:52:06.350617",
    "2021-10-05 10:52:06.350620",
    "2021-10-05 10:52:06.350623",
    "2021-10-05 10:52:06.350626"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
